# MS: vyplneni a kontrola dalsiho formulare
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update selection on the existing "NewPerson1" sheet
$ws1.Range("C2").Select()

# Refresh the header fill to the Accent4 theme color (same orange, now theme-linked)
$ws1.Range("A1:G1").Interior.ThemeColor = 8

# Add the new "CZ_HPP" worksheet right after "NewPerson1"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "CZ_HPP"

# --- Header row (row 1) ---
$ws2.Range("A1").Value = "Titul"
$ws2.Range("B1").Value = "Jméno"
$ws2.Range("C1").Value = "Příjmení"
$ws2.Range("D1").Value = "Rodné číslo"
$ws2.Range("E1").Value = "Číslo OP"
$ws2.Range("F1").Value = "E-mail"
$ws2.Range("G1").Value = "Telefon"
$ws2.Range("H1").Value = "Ulice"
$ws2.Range("I1").Value = "PSČ"
$ws2.Range("J1").Value = "Město"

# Orange (Accent 4 theme) fill on the header row, same accent color used on NewPerson1
$ws2.Range("A1:J1").Interior.ThemeColor = 8

# --- Data row (row 2) ---
$ws2.Range("A2").Value = "Ing."

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Formula = "=NewPerson1!A2"

$ws2.Range("C2").NumberFormat = "@"
$ws2.Range("C2").Formula = "=NewPerson1!B2"

$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").Formula = "=NewPerson1!C2"

$ws2.Range("G2").NumberFormat = "@"
$ws2.Range("G2").Formula = "=NewPerson1!D2"

$ws2.Range("H2").Value = "Testerská 1234"
$ws2.Range("I2").Value = 33333
$ws2.Range("J2").Value = "Praha"

$ws2.Range("D2").Value = "001017/6573"

$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "222222222"

# CZ_HPP becomes the active sheet/tab, with E3 selected
$ws2.Range("E3").Select()
